# This script applies the "face -> book" stimulus rename and the
# correct_ans single-letter -> full-word expansion (r->right, y->left,
# b->center) described in the commit's accompanying diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the extent of the data (header row 1 + data rows).
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Mapping for the abbreviated "correct_ans" codes in column L.
$ansMap = @{ "r" = "right"; "y" = "left"; "b" = "center" }

# Columns that may contain the "face//face_NN.jpg" style stimulus paths.
$fileCols = @("A", "B", "C", "D")

for ($row = 2; $row -le $lastRow; $row++) {

    foreach ($col in $fileCols) {
        $cell = $ws.Range($col + $row)
        $val = $cell.Value()
        if ($val -ne $null) {
            $text = $val.ToString()
            if ($text.Contains("face")) {
                $cell.Value = $text.Replace("face", "book")
            }
        }
    }

    $lCell = $ws.Range("L" + $row)
    $lVal = $lCell.Value()
    if ($lVal -ne $null) {
        $lText = $lVal.ToString()
        if ($ansMap.ContainsKey($lText)) {
            $lCell.Value = $ansMap[$lText]
        }
    }
}
